$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.585.66"
$ws.Range("E2").Value = "  +6.96%  "
$ws.Range("D3").Value = "2.582.68"
$ws.Range("E3").Value = "  +8.67%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'504.81"
$ws.Range("E5").Value = "  +5.64%  "
$ws.Range("D6").Value = "'157.34"
$ws.Range("E6").Value = "  +6.95%  "
$ws.Range("E7").Value = "  +24.01%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D9").Value = "2.580.41"
$ws.Range("E9").Value = "  +8.41%  "
$ws.Range("D10").Value = "'6.19"
$ws.Range("E10").Value = "  +14.04%  "
$ws.Range("E11").Value = "  +5.97%  "
$ws.Range("D12").Value = "'0.341"
$ws.Range("E12").Value = "  +5.59%  "
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "3.028.87"
$ws.Range("E14").Value = "  +8.63%  "
$ws.Range("D15").Value = "59.422.16"
$ws.Range("D16").Value = "'21.86"
$ws.Range("E17").Value = "  +4.17%  "
$ws.Range("D18").Value = "2.580.98"
$ws.Range("E18").Value = "  +8.26%  "
$ws.Range("D19").Value = "'4.73"
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("D20").Value = "'334.46"
$ws.Range("E20").Value = "  +6.15%  "
$ws.Range("D21").Value = "'10.38"
$ws.Range("E21").Value = "  +7.25%  "
$ws.Range("D22").Value = "'6.07"
$ws.Range("E22").Value = "  +7.32%  "
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").Value = "'60.31"
$ws.Range("E24").Value = "  +6.10%  "
$ws.Range("D25").Value = "'0.415"
$ws.Range("E25").Value = "  +4.97%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  +7.06%  "
$ws.Range("D27").Value = "2.685.64"
$ws.Range("E27").Value = "  +8.19%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'7.46"
$ws.Range("E29").Value = "  +3.02%  "
$ws.Range("D30").Value = "0.0₃0820"
$ws.Range("E30").Value = "  +6.51%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "'156.75"
$ws.Range("E32").Value = "  +6.03%  "
$ws.Range("D33").Value = "'19.43"
$ws.Range("E33").Value = "  +7.91%  "
$ws.Range("E34").Value = "  +5.37%  "
$ws.Range("E35").Value = "  +8.18%  "
$ws.Range("D36").Value = "'3.95"
$ws.Range("E36").Value = "  +10.64%  "
$ws.Range("E37").Value = "  +8.28%  "
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").Value = "'3.78"
$ws.Range("E39").Value = "  +11.71%  "
$ws.Range("E40").Value = "  +7.49%  "
$ws.Range("D41").Value = "'35.03"
$ws.Range("E41").Value = "  +4.85%  "
$ws.Range("D42").Value = "'291.83"
$ws.Range("E42").Value = "  +14.08%  "
$ws.Range("E43").Value = "  +7.58%  "
$ws.Range("E44").Value = "  +7.30%  "
$ws.Range("E45").Value = "  +4.66%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'19.29"
$ws.Range("E47").Value = "  +14.19%  "
$ws.Range("D48").Value = "'0.0238"
$ws.Range("E48").Value = "  +6.70%  "
$ws.Range("D49").Value = "'4.77"
$ws.Range("E49").Value = "  +5.23%  "
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "'0.711"
$ws.Range("E51").Value = "  +11.47%  "
